# Jun's file updates for all IO data and others
#
# Replaces the hard-coded placeholder "1" values in the PoFDCtAE sheet with
# live formulas that pull the percentage of fuel-demand changes that affect
# exports from the "Data from BFPIaE" sheet (columns B/D, and D/(D+E) where
# a combined share is needed). Dependent cells (the "1 - x" helper column R)
# recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PoFDCtAE")

$ws.Range("C3").Formula  = "='Data from BFPIaE'!D5/'Data from BFPIaE'!B5"
$ws.Range("D4").Formula  = "='Data from BFPIaE'!D6/'Data from BFPIaE'!B6"
$ws.Range("I9").Formula  = "='Data from BFPIaE'!D11/'Data from BFPIaE'!B11"
$ws.Range("J10").Formula = "='Data from BFPIaE'!D12/SUM('Data from BFPIaE'!D12:E12)"
$ws.Range("K11").Formula = "='Data from BFPIaE'!D13/SUM('Data from BFPIaE'!D13:E13)"
$ws.Range("L12").Formula = "='Data from BFPIaE'!D14/'Data from BFPIaE'!B14"
$ws.Range("M13").Formula = "='Data from BFPIaE'!D15/'Data from BFPIaE'!B15"
$ws.Range("N14").Formula = "='Data from BFPIaE'!D16/SUM('Data from BFPIaE'!D16:E16)"
$ws.Range("S19").Formula = "='Data from BFPIaE'!D21/SUM('Data from BFPIaE'!D21:E21)"
$ws.Range("T20").Formula = "='Data from BFPIaE'!D22/SUM('Data from BFPIaE'!D22:E22)"

$excel.Calculate()
